$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values
$ws.Range("B2").Value = 95.695904314898485
$ws.Range("C2").Value = 93.967282776995887
$ws.Range("D2").Value = 93.849075396126196
$ws.Range("E2").Value = 94.821777757501053

# Row 3 values
$ws.Range("B3").Value = 94.174636240268143
$ws.Range("C3").Value = 94.067146827857201
$ws.Range("D3").Value = 92.153541463218403
$ws.Range("E3").Value = 96.43580666409791

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
